$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 73) with the latest quotation data (2025-11-16)
$ws.Range("A73").Value = 45977
$ws.Range("B73").Value = "22,2279"
$ws.Range("C73").Value = "15,8858"
$ws.Range("D73").Value = "15,6322"
$ws.Range("E73").Value = "15,6322"

# Match the date-formatted style used by the rest of column A
$ws.Range("A2").Copy()
$ws.Range("A73").PasteSpecial(-4122)
